$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.568.86"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.808.91"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.70"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.596"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "37.49"
$ws.Range("E8").Value = "  +7.29%  "
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0681"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "2.070.69"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.30"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "1.819.41"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "34.533.21"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.67"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.59"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "0.0₃0775"
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("E24").Value = "  +4.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.06"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.83"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.30"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.120"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.92"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "1.363.91"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("E36").Value = "  -4.68%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  -5.63%  "
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "80.79"
$ws.Range("E42").Value = "  -3.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.938"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("E44").Value = "  +5.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.69"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").Value = "1.970.70"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.62"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  -4.70%  "
